# Generate Report for Handoff
# Update status from "In Translation" to "Ready for handoff" and refresh
# the handoff timestamps across the Overview, zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet (row 2): Status columns B2/C2, Latest Handoff Datetime D2
$wsOverview.Range("B2").Value = "Ready for handoff"
$wsOverview.Range("C2").Value = "Ready for handoff"
$wsOverview.Range("D2").Value = "2016-30-20 12:30:07"

# zh-cn sheet (row 2): Status C2, Latest Handoff Datetime E2
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("E2").Value = "2016-03-20 12:30:01"

# de-de sheet (row 2): Status C2, Latest Handoff Datetime E2
$wsDeDe.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("E2").Value = "2016-03-20 12:30:07"
